$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-29 (Gender column C) -> Male
$ws.Range("C2:C29").Value = "Male"

# Rows 30-57 (Gender column C) -> Female
$ws.Range("C30:C57").Value = "Female"
